$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new formula cell C1 (text formula producing "-Wolfram") — one
# label per .csv / process, per the commit message.
$ws.Range("C1").Formula = '="-Wolfram"'

# Move the selection/active cell to where the author's session ended up
# (this also grows <dimension> to A1:C1 and <row spans> to 1:3).
$ws.Range("C11").Select()
